$d = $word.ActiveDocument

# --- Paragraph 1: "This is a Microsoft word document." ---
$p1 = $d.Paragraphs(1)
$r = $p1.Range

# 1) two trailing spaces appended to the original (unformatted) run's text
$r.InsertAfter("  ")

# Insert the new blank paragraph right now, while formatting is still the
# plain/default formatting -- this keeps the new paragraph mark's run free
# of the red color that the "(This is a change...)" text below will use.
$r.InsertParagraphAfter()

# Re-acquire paragraph 1 (now holding just the original sentence + 2 spaces)
# and append the three red runs to it, before its own paragraph mark.
$p1 = $d.Paragraphs(1)
$r = $p1.Range

# 2) red run: "(This is a change " + en-dash + " Version for branch "
$start = $r.End - 1
$r.InsertAfter("(This is a change " + [char]0x2013 + " Version for branch ")
$end = $r.End - 1
$d.Range($start, $end).Font.Color = 192

# 3) red run: "main"
$start = $r.End - 1
$r.InsertAfter("main")
$end = $r.End - 1
$d.Range($start, $end).Font.Color = 192

# 4) red run: ")"
$start = $r.End - 1
$r.InsertAfter(")")
$end = $r.End - 1
$d.Range($start, $end).Font.Color = 192
